# auto: removing some labels from the patient card
#
# The "survey" sheet has two "note" rows that only existed to render
# extra labels on the patient card:
#   - row 40: NO_LABEL "nick"    -> "Nickname: **${aka_ctx}**"
#   - row 42: NO_LABEL "gender_n" -> "Gender Identity: **${gender_ctx}**"
#
# Both rows are removed entirely (including their shared strings),
# which shifts every row below them up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Delete the lower row first so row 40's index isn't invalidated by the
# first deletion.
$ws.Rows.Item(42).EntireRow.Delete()
$ws.Rows.Item(40).EntireRow.Delete()
